# ClueBoardLayout.xlsx edit
# - Row 30 on Sheet1 now holds a 0-based column index sequence starting in
#   column A (previously it started one column to the right, in column B,
#   and only ran 0-19; now it runs A:U with values 0-20).
# - The active sheet view's scroll position / selection moved: the window
#   now scrolls so column D is left-most, and the active cell/selection is
#   T33 (previously topLeftCell A2 / selection Y28).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 30: shift the 0..19 sequence one column right (A:U => 0..20) ---
$row30Values = @(0,1,2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20)
for ($i = 0; $i -lt $row30Values.Length; $i++) {
    $ws.Cells.Item(30, $i + 1).Value = $row30Values[$i]
}

# --- Sheet view: scroll so column D is the left-most visible column, and
#     move the selection to T33 ---
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 4
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("T33").Select()
